$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 74356.71000000001
$ws.Range("I28").Value = 369.8
$ws.Range("J28").Value = 259324
$ws.Range("K28").Value = 369.8
$ws.Range("L28").Value = 259324
$ws.Range("M28").Value = 115.2
$ws.Range("N28").Value = -260294

$ws.Range("H33").Value = 37342.07
$ws.Range("I33").Value = 45935.547
$ws.Range("K33").Value = 45935.547
$ws.Range("M33").Value = -45706.547

$ws.Range("H39").Value = 297.5238
$ws.Range("J39").Value = 595.1111
$ws.Range("L39").Value = 1785.3333
$ws.Range("N39").Value = -2377.3333

$ws.Range("H98").Value = 559.63635
$ws.Range("J98").Value = 999.5
$ws.Range("L98").Value = 999.5
$ws.Range("N98").Value = -3995.5

$ws.Range("H115").Value = 627.6
$ws.Range("I115").Value = 627.6
$ws.Range("K115").Value = 1882.8
$ws.Range("M115").Value = -315.8000000000002

$ws.Range("H122").Value = 559.63635
$ws.Range("J122").Value = 999.5
$ws.Range("L122").Value = 2998.5
$ws.Range("N122").Value = -7898.5

$ws.Range("H125").Value = 86106.234
$ws.Range("I125").Value = 8085.25
$ws.Range("K125").Value = 72767.25
$ws.Range("M125").Value = -70307.25

$ws.Range("H137").Value = 2717.1
$ws.Range("I137").Value = 1790.2354
$ws.Range("K137").Value = 5370.706200000001
$ws.Range("M137").Value = -2820.706200000001

$ws.Range("H138").Value = 3005.432
$ws.Range("I138").Value = 2116.4119
$ws.Range("J138").Value = 6028.1
$ws.Range("K138").Value = 6349.2357
$ws.Range("L138").Value = 18084.3
$ws.Range("M138").Value = -1209.2357
$ws.Range("N138").Value = -28364.3

$ws.Range("H141").Value = 3661.3333
$ws.Range("I141").Value = 3493.6667
$ws.Range("J141").Value = 4499.6665
$ws.Range("K141").Value = 10481.0001
$ws.Range("L141").Value = 13498.9995
$ws.Range("M141").Value = -5301.000100000001
$ws.Range("N141").Value = -23858.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6738.6665
$ws.Range("I32").Value = 2878.8728
$ws.Range("K32").Value = 2878.8728
$ws.Range("M32").Value = -2591.8728

$ws.Range("H61").Value = 6442.3794
$ws.Range("I61").Value = 5616.2915
$ws.Range("K61").Value = 5616.2915
$ws.Range("M61").Value = -5404.2915

$ws.Range("H74").Value = 3230.0625
$ws.Range("I74").Value = 1654.12
$ws.Range("K74").Value = 1654.12
$ws.Range("M74").Value = -780.1199999999999

$ws.Range("H77").Value = 3230.0625
$ws.Range("I77").Value = 1654.12
$ws.Range("K77").Value = 8270.599999999999
$ws.Range("M77").Value = -3902.599999999999

$ws.Range("H122").Value = 3530.8125
$ws.Range("I122").Value = 3373.8462
$ws.Range("K122").Value = 10121.5386
$ws.Range("M122").Value = -7671.5386

$ws.Range("H132").Value = 3212.4119
$ws.Range("I132").Value = 1799.25
$ws.Range("K132").Value = 5397.75
$ws.Range("M132").Value = -2867.75

$ws.Range("H136").Value = 6442.3794
$ws.Range("I136").Value = 5616.2915
$ws.Range("K136").Value = 16848.8745
$ws.Range("M136").Value = -14298.8745

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H134").Value = 2514.125
$ws.Range("I134").Value = 1856.9736
$ws.Range("K134").Value = 5570.9208
$ws.Range("M134").Value = -3035.9208

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 681.05884
$ws.Range("I7").Value = 745.3043
$ws.Range("J7").Value = 546.7273
$ws.Range("K7").Value = 745.3043
$ws.Range("L7").Value = 546.7273
$ws.Range("M7").Value = -632.3043
$ws.Range("N7").Value = -772.7273

$ws.Range("H31").Value = 7960.1914
$ws.Range("I31").Value = 2847.4
$ws.Range("J31").Value = 16982.766
$ws.Range("K31").Value = 2847.4
$ws.Range("L31").Value = 16982.766
$ws.Range("M31").Value = -2552.4
$ws.Range("N31").Value = -17572.766

$ws.Range("H34").Value = 7960.1914
$ws.Range("I34").Value = 2847.4
$ws.Range("J34").Value = 16982.766
$ws.Range("K34").Value = 2847.4
$ws.Range("L34").Value = 16982.766
$ws.Range("M34").Value = -2645.4
$ws.Range("N34").Value = -17386.766

$ws.Range("H58").Value = 5879.2144
$ws.Range("I58").Value = 2591.524
$ws.Range("J58").Value = 15742.286
$ws.Range("K58").Value = 2591.524
$ws.Range("L58").Value = 15742.286
$ws.Range("M58").Value = -2388.524
$ws.Range("N58").Value = -16148.286

$ws.Range("H103").Value = 124192
$ws.Range("I103").Value = 15993
$ws.Range("J103").Value = 178291.5
$ws.Range("K103").Value = 15993
$ws.Range("L103").Value = 178291.5
$ws.Range("M103").Value = -14821
$ws.Range("N103").Value = -180635.5

$ws.Range("H107").Value = 1287.9412
$ws.Range("I107").Value = 1073.0714
$ws.Range("J107").Value = 2290.6667
$ws.Range("K107").Value = 1073.0714
$ws.Range("L107").Value = 2290.6667
$ws.Range("M107").Value = 846.9286
$ws.Range("N107").Value = -6130.6667

$ws.Range("H108").Value = 199950
$ws.Range("J108").Value = 199950
$ws.Range("L108").Value = 199950
$ws.Range("N108").Value = -207630

$ws.Range("H132").Value = 4281.757
$ws.Range("I132").Value = 3045.5862
$ws.Range("K132").Value = 9136.758600000001
$ws.Range("M132").Value = -6606.758600000001

$ws.Range("H134").Value = 6035.5713
$ws.Range("I134").Value = 2899.111
$ws.Range("K134").Value = 8697.332999999999
$ws.Range("M134").Value = -6162.332999999999

$ws.Range("H136").Value = 5879.2144
$ws.Range("I136").Value = 2591.524
$ws.Range("J136").Value = 15742.286
$ws.Range("K136").Value = 7774.572
$ws.Range("L136").Value = 47226.858
$ws.Range("M136").Value = -5224.572
$ws.Range("N136").Value = -52326.858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 421.72415
$ws.Range("I2").Value = 32.105263
$ws.Range("K2").Value = 192.631578
$ws.Range("M2").Value = -79.63157799999999

$ws.Range("H4").Value = 62501120
$ws.Range("J4").Value = 1640.4
$ws.Range("L4").Value = 4921.200000000001
$ws.Range("N4").Value = -5145.200000000001

$ws.Range("H70").Value = 12115.733
$ws.Range("I70").Value = 9248.727999999999
$ws.Range("J70").Value = 20000
$ws.Range("K70").Value = 27746.184
$ws.Range("L70").Value = 60000
$ws.Range("M70").Value = -27431.184
$ws.Range("N70").Value = -60630

$ws.Range("H73").Value = 12115.733
$ws.Range("I73").Value = 9248.727999999999
$ws.Range("J73").Value = 20000
$ws.Range("K73").Value = 27746.184
$ws.Range("L73").Value = 60000
$ws.Range("M73").Value = -26654.184
$ws.Range("N73").Value = -62184

$ws.Range("H75").Value = 6039.8
$ws.Range("I75").Value = 699
$ws.Range("K75").Value = 2097
$ws.Range("M75").Value = -1099

$ws.Range("H78").Value = 6039.8
$ws.Range("I78").Value = 699
$ws.Range("K78").Value = 6291
$ws.Range("M78").Value = -1299

$ws.Range("H111").Value = 5000
$ws.Range("I111").Value = 5000
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 15000
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H131").Value = 729126.4
$ws.Range("I131").Value = 970.1667
$ws.Range("K131").Value = 2910.5001
$ws.Range("M131").Value = 2129.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3443.4285
$ws.Range("I102").Value = 1998.8182
$ws.Range("K102").Value = 1998.8182
$ws.Range("M102").Value = -376.8181999999999

$ws.Range("H122").Value = 10839.8
$ws.Range("I122").Value = 9542.571
$ws.Range("K122").Value = 28627.713
$ws.Range("M122").Value = -26177.713

$ws.Range("H132").Value = 5884.1377
$ws.Range("I132").Value = 5591.88
$ws.Range("J132").Value = 7710.75
$ws.Range("K132").Value = 16775.64
$ws.Range("L132").Value = 23132.25
$ws.Range("M132").Value = -14245.64
$ws.Range("N132").Value = -28192.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3800
$ws.Range("J46").Value = 4825
$ws.Range("L46").Value = 4825
$ws.Range("N46").Value = -5201

$ws.Range("H109").Value = 154250
$ws.Range("I109").Value = 30000
$ws.Range("J109").Value = 174958.33
$ws.Range("K109").Value = 30000
$ws.Range("L109").Value = 174958.33
$ws.Range("M109").Value = -28613
$ws.Range("N109").Value = -177732.33

$ws.Range("H122").Value = 6248.2144
$ws.Range("I122").Value = 5734.375
$ws.Range("K122").Value = 17203.125
$ws.Range("M122").Value = -14753.125

$ws.Range("H132").Value = 6369.476
$ws.Range("J132").Value = 11222
$ws.Range("L132").Value = 33666
$ws.Range("N132").Value = -38726
